$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$f2 = @'
<rpc-reply message-id="urn:uuid:4ce91b62-500b-4dfb-bf13-8b5f76326c86">
  <data/>
</rpc-reply>

'@

$g2 = @'
  <edit-config>
    <target>
     <candidate/>
    </target>
    <config>
      <network-instances xmlns="http://openconfig.net/yang/network-instance">
        <network-instance>
          <name>Prueba_LxVPN</name>
          <config>
            <name>Prueba_LxVPN</name>
            <type xmlns:oc-ni-types="http://openconfig.net/yang/network-instance-types">oc-ni-types:L3VRF</type>
          </config>
          <protocols>
            <protocol>
              <identifier xmlns:oc-pol-types="http://openconfig.net/yang/policy-types">oc-pol-types:BGP</identifier>
              <name>default</name>
              <config>
                <identifier xmlns:oc-pol-types="http://openconfig.net/yang/policy-types">oc-pol-types:BGP</identifier>
                <name>default</name>
              </config>
              <bgp>
                <global>
                  <config>
                    <as>65000</as>
                  </config>
                </global>
                <neighbors>
                  <neighbor>
                    <neighbor-address>192.168.1.2</neighbor-address>
                    <config>
                      <neighbor-address>192.168.1.2</neighbor-address>
                      <peer-as>65123</peer-as>
                    </config>
                    <ebgp-multihop>
                      <config>
                        <multihop-ttl>3</multihop-ttl>
                      </config>
                    </ebgp-multihop>
                  </neighbor>
                </neighbors>
              </bgp>
            </protocol>
          </protocols>
        </network-instance>
      </network-instances>
    </config>
  </edit-config>
'@

$ws.Range("F2").Value = $f2
$ws.Range("G2").Value = $g2
